$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 621, pushing the former rows 621-662 down to 622-663.
$ws.Rows(621).Insert()

# Populate the newly inserted row 621 with the new weekly record.
$ws.Cells.Item(621, 1).Value  = 9
$ws.Cells.Item(621, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(621, 3).Value  = "Metropolitana"
$ws.Cells.Item(621, 4).Value  = 45265
$ws.Cells.Item(621, 5).Value  = 13
$ws.Cells.Item(621, 6).Value  = 100112039
$ws.Cells.Item(621, 7).Value  = "Ciboulette"
$ws.Cells.Item(621, 8).Value  = "Sin especificar"
$ws.Cells.Item(621, 9).Value  = "Primera"
$ws.Cells.Item(621, 10).Value = 430
$ws.Cells.Item(621, 11).Value = 1000
$ws.Cells.Item(621, 12).Value = 1200
$ws.Cells.Item(621, 13).Value = 1100
$ws.Cells.Item(621, 14).Value = "$/docena de atados"
$ws.Cells.Item(621, 15).Value = "Región Metropolitana"
$ws.Cells.Item(621, 16).Value = 367
$ws.Cells.Item(621, 17).Value = 3
$ws.Cells.Item(621, 18).Value = "Hortaliza"
